$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.526.38"
$ws.Range("E2").Value = "  +4.19%  "

# Row 3
$ws.Range("D3").Value = "1.791.93"
$ws.Range("E3").Value = "  +0.70%  "

# Row 4
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.49%  "

# Row 5
$ws.Range("D5").Value = "'313.35"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.37%  "

# Row 7
$ws.Range("D7").Value = "'0.5339"
$ws.Range("E7").Value = "  +0.44%  "

# Row 8
$ws.Range("D8").Value = "'0.3793"
$ws.Range("E8").Value = "  +0.84%  "

# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.07509"
$ws.Range("E9").Value = "  +1.34%  "

# Row 10
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'42.49"
$ws.Range("E10").Value = "  -0.81%  "

# Row 11
$ws.Range("D11").Value = "'1.116"
$ws.Range("E11").Value = "  +2.00%  "

# Row 12
$ws.Range("D12").Value = "'1.004"
$ws.Range("E12").Value = "  +0.32%  "

# Row 13
$ws.Range("D13").Value = "'21.03"
$ws.Range("E13").Value = "  +1.64%  "

# Row 14
$ws.Range("D14").Value = "'6.161"
$ws.Range("E14").Value = "  +0.85%  "

# Row 15
$ws.Range("D15").Value = "'7.382"
$ws.Range("E15").Value = "  +5.61%  "

# Row 16
$ws.Range("D16").Value = "1.795.73"
$ws.Range("E16").Value = "  +0.77%  "

# Row 17
$ws.Range("D17").Value = "'90.29"
$ws.Range("E17").Value = "  +0.58%  "

# Row 18
$ws.Range("D18").Value = "'0.00001065"
$ws.Range("E18").Value = "  +0.83%  "

# Row 19
$ws.Range("D19").Value = "'0.06436"
$ws.Range("E19").Value = "  +0.02%  "

# Row 20
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.37%  "

# Row 21
$ws.Range("D21").Value = "'17.27"
$ws.Range("E21").Value = "  +2.82%  "

# Row 22
$ws.Range("D22").Value = "'5.923"
$ws.Range("E22").Value = "  +0.26%  "

# Row 23
$ws.Range("D23").Value = "28.606.55"
$ws.Range("E23").Value = "  +4.30%  "

# Row 24
$ws.Range("D24").Value = "'11.22"
$ws.Range("E24").Value = "  +0.17%  "

# Row 25
$ws.Range("D25").Value = "'2.135"
$ws.Range("E25").Value = "  +1.85%  "

# Row 26
$ws.Range("D26").Value = "'160.58"
$ws.Range("E26").Value = "  +3.23%  "

# Row 27
$ws.Range("D27").Value = "'20.45"
$ws.Range("E27").Value = "  +1.11%  "

# Row 28
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.002.75"
$ws.Range("E28").Value = "  +0.73%  "

# Row 29
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.356"
$ws.Range("E29").Value = "  -0.58%  "

# Row 30
$ws.Range("D30").Value = "'123.19"
$ws.Range("E30").Value = "  +1.42%  "

# Row 31
$ws.Range("D31").Value = "'1.124"
$ws.Range("E31").Value = "  +3.39%  "

# Row 32
$ws.Range("D32").Value = "'0.1015"
$ws.Range("E32").Value = "  -1.33%  "

# Row 33
$ws.Range("D33").Value = "'5.664"
$ws.Range("E33").Value = "  +1.29%  "

# Row 34
$ws.Range("D34").Value = "'3.656"
$ws.Range("E34").Value = "  +0.58%  "

# Row 35
$ws.Range("D35").Value = "'0.2306"
$ws.Range("E35").Value = "  +12.40%  "

# Row 36
$ws.Range("D36").Value = "'0.06560"
$ws.Range("E36").Value = "  +9.72%  "

# Row 37
$ws.Range("D37").Value = "'0.02321"
$ws.Range("E37").Value = "  +2.83%  "

# Row 38
$ws.Range("D38").Value = "'8.695"
$ws.Range("E38").Value = "  +5.59%  "

# Row 39
$ws.Range("D39").Value = "'5.088"
$ws.Range("E39").Value = "  +3.44%  "

# Row 40
$ws.Range("D40").Value = "'11.46"
$ws.Range("E40").Value = "  +1.53%  "

# Row 41
$ws.Range("D41").Value = "'0.6314"
$ws.Range("E41").Value = "  +2.96%  "

# Row 42
$ws.Range("D42").Value = "'1.211"
$ws.Range("E42").Value = "  +6.66%  "

# Row 43
$ws.Range("D43").Value = "'1.004"
$ws.Range("E43").Value = "  +0.41%  "

# Row 44
$ws.Range("D44").Value = "'1.393"
$ws.Range("E44").Value = "  -2.88%  "

# Row 45
$ws.Range("D45").Value = "'13.44"
$ws.Range("E45").Value = "  +1.64%  "

# Row 46
$ws.Range("D46").Value = "'0.5926"
$ws.Range("E46").Value = "  +2.27%  "

# Row 47
$ws.Range("D47").Value = "'3.670"
$ws.Range("E47").Value = "  +1.19%  "

# Row 48
$ws.Range("D48").Value = "'124.65"
$ws.Range("E48").Value = "  +2.58%  "

# Row 49
$ws.Range("D49").Value = "'1.979"
$ws.Range("E49").Value = "  +4.38%  "

# Row 50
$ws.Range("D50").Value = "'1.154"
$ws.Range("E50").Value = "  +2.96%  "

# Row 51
$ws.Range("D51").Value = "'0.06919"
$ws.Range("E51").Value = "  +2.83%  "

Write-Host "Applied crypto list update"
